# Applies the cryptos.xlsx data-refresh diff: updated Price/Volume(1h) figures
# for every coin row, plus a re-rank of 5 coins (rows 38-42) whose Coin/Link/
# Price/Volume cells rotate to a new row position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.436.97'; AsText = $true },
    @{ Cell = 'E2'; Value = '  +0.14%  '; AsText = $false },
    @{ Cell = 'D3'; Value = '1.636.15'; AsText = $true },
    @{ Cell = 'E3'; Value = '  -0.90%  '; AsText = $false },
    @{ Cell = 'E4'; Value = '  +0.01%  '; AsText = $false },
    @{ Cell = 'D5'; Value = '212.53'; AsText = $true },
    @{ Cell = 'E5'; Value = '  -0.61%  '; AsText = $false },
    @{ Cell = 'D6'; Value = '0.532'; AsText = $true },
    @{ Cell = 'E6'; Value = '  +4.55%  '; AsText = $false },
    @{ Cell = 'E7'; Value = '  +0.03%  '; AsText = $false },
    @{ Cell = 'D8'; Value = '23.01'; AsText = $true },
    @{ Cell = 'E8'; Value = '  -5.05%  '; AsText = $false },
    @{ Cell = 'D9'; Value = '0.257'; AsText = $true },
    @{ Cell = 'E9'; Value = '  -2.32%  '; AsText = $false },
    @{ Cell = 'D10'; Value = '0.0610'; AsText = $true },
    @{ Cell = 'E10'; Value = '  -0.86%  '; AsText = $false },
    @{ Cell = 'E11'; Value = '  +1.09%  '; AsText = $false },
    @{ Cell = 'D12'; Value = '1.867.30'; AsText = $true },
    @{ Cell = 'E12'; Value = '  -0.89%  '; AsText = $false },
    @{ Cell = 'D13'; Value = '1.646.44'; AsText = $true },
    @{ Cell = 'E13'; Value = '  -0.33%  '; AsText = $false },
    @{ Cell = 'D14'; Value = '0.581'; AsText = $true },
    @{ Cell = 'E14'; Value = '  +3.28%  '; AsText = $false },
    @{ Cell = 'D15'; Value = '4.02'; AsText = $true },
    @{ Cell = 'E15'; Value = '  -1.77%  '; AsText = $false },
    @{ Cell = 'D16'; Value = '64.14'; AsText = $true },
    @{ Cell = 'E16'; Value = '  -2.36%  '; AsText = $false },
    @{ Cell = 'D17'; Value = '27.429.81'; AsText = $true },
    @{ Cell = 'E17'; Value = '  +0.15%  '; AsText = $false },
    @{ Cell = 'D18'; Value = '229.03'; AsText = $true },
    @{ Cell = 'E18'; Value = '  -2.50%  '; AsText = $false },
    @{ Cell = 'D19'; Value = '0.0₃0723'; AsText = $true },
    @{ Cell = 'E19'; Value = '  -0.44%  '; AsText = $false },
    @{ Cell = 'D20'; Value = '7.60'; AsText = $true },
    @{ Cell = 'E20'; Value = '  +1.04%  '; AsText = $false },
    @{ Cell = 'E21'; Value = '  -0.04%  '; AsText = $false },
    @{ Cell = 'D22'; Value = '4.31'; AsText = $true },
    @{ Cell = 'E22'; Value = '  -2.41%  '; AsText = $false },
    @{ Cell = 'D23'; Value = '9.74'; AsText = $true },
    @{ Cell = 'E23'; Value = '  +5.66%  '; AsText = $false },
    @{ Cell = 'D24'; Value = '1.96'; AsText = $true },
    @{ Cell = 'E24'; Value = '  -3.37%  '; AsText = $false },
    @{ Cell = 'D25'; Value = '149.42'; AsText = $true },
    @{ Cell = 'E25'; Value = '  +2.42%  '; AsText = $false },
    @{ Cell = 'D26'; Value = '7.00'; AsText = $true },
    @{ Cell = 'E26'; Value = '  -2.33%  '; AsText = $false },
    @{ Cell = 'E27'; Value = '  +1.72%  '; AsText = $false },
    @{ Cell = 'E28'; Value = '  -0.07%  '; AsText = $false },
    @{ Cell = 'D29'; Value = '15.55'; AsText = $true },
    @{ Cell = 'E29'; Value = '  -3.23%  '; AsText = $false },
    @{ Cell = 'E30'; Value = '  -0.69%  '; AsText = $false },
    @{ Cell = 'D31'; Value = '0.0488'; AsText = $true },
    @{ Cell = 'E31'; Value = '  -1.94%  '; AsText = $false },
    @{ Cell = 'D32'; Value = '3.29'; AsText = $true },
    @{ Cell = 'E32'; Value = '  -0.44%  '; AsText = $false },
    @{ Cell = 'E33'; Value = '  +2.83%  '; AsText = $false },
    @{ Cell = 'D34'; Value = '1.423.98'; AsText = $true },
    @{ Cell = 'E34'; Value = '  -2.40%  '; AsText = $false },
    @{ Cell = 'D35'; Value = '1.59'; AsText = $true },
    @{ Cell = 'E35'; Value = '  +2.44%  '; AsText = $false },
    @{ Cell = 'E36'; Value = '  -1.89%  '; AsText = $false },
    @{ Cell = 'D37'; Value = '0.570'; AsText = $true },
    @{ Cell = 'E37'; Value = '  -0.05%  '; AsText = $false },
    @{ Cell = 'B38'; Value = 'VeChain'; AsText = $false },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; AsText = $false },
    @{ Cell = 'D38'; Value = '0.0168'; AsText = $true },
    @{ Cell = 'E38'; Value = '  -1.17%  '; AsText = $false },
    @{ Cell = 'B39'; Value = 'ARBITRUM'; AsText = $false },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; AsText = $false },
    @{ Cell = 'D39'; Value = '0.874'; AsText = $true },
    @{ Cell = 'E39'; Value = '  -4.19%  '; AsText = $false },
    @{ Cell = 'B40'; Value = 'TrustWalletToken'; AsText = $false },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; AsText = $false },
    @{ Cell = 'D40'; Value = '0.882'; AsText = $true },
    @{ Cell = 'E40'; Value = '  +12.09%  '; AsText = $false },
    @{ Cell = 'B41'; Value = 'WEMIXToken'; AsText = $false },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; AsText = $false },
    @{ Cell = 'D41'; Value = '1.03'; AsText = $true },
    @{ Cell = 'E41'; Value = '  -1.31%  '; AsText = $false },
    @{ Cell = 'B42'; Value = 'PaxDollar'; AsText = $false },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; AsText = $false },
    @{ Cell = 'D42'; Value = '1.00'; AsText = $true },
    @{ Cell = 'E42'; Value = '  +0.03%  '; AsText = $false },
    @{ Cell = 'D43'; Value = '5.52'; AsText = $true },
    @{ Cell = 'E43'; Value = '  +1.52%  '; AsText = $false },
    @{ Cell = 'E44'; Value = '  +0.46%  '; AsText = $false },
    @{ Cell = 'D45'; Value = '64.84'; AsText = $true },
    @{ Cell = 'E45'; Value = '  -0.54%  '; AsText = $false },
    @{ Cell = 'D46'; Value = '1.777.51'; AsText = $true },
    @{ Cell = 'E46'; Value = '  -0.71%  '; AsText = $false },
    @{ Cell = 'D47'; Value = '1.66'; AsText = $true },
    @{ Cell = 'E47'; Value = '  -3.24%  '; AsText = $false },
    @{ Cell = 'D48'; Value = '85.72'; AsText = $true },
    @{ Cell = 'E48'; Value = '  -2.99%  '; AsText = $false },
    @{ Cell = 'E49'; Value = '  -0.17%  '; AsText = $false },
    @{ Cell = 'D50'; Value = '0.0990'; AsText = $true },
    @{ Cell = 'E50'; Value = '  -1.85%  '; AsText = $false },
    @{ Cell = 'E51'; Value = '  -1.10%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
